# Applies the "multiple inheritance via arrows" diagram restructuring.
# Strategy:
#  1. Capture stable format templates (cellXf style indices) for the two
#     existing left-border accent styles used throughout the sheet
#     (the "thick black left border" style and the "thin red left
#     border" style), copying them to scratch cells far outside the
#     used range so row/column moves below don't disturb them.
#  2. Wipe the existing data block (rows 2-17, columns A-F).
#  3. Re-create every cell of the new layout (rows 2-24) with the
#     correct text and formatting, applying the saved format templates
#     via PasteSpecial (formats only) so the existing style indices in
#     styles.xml are reused rather than duplicated.
#  4. Clean up the scratch template cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Save format templates from existing styled cells -----------------
$tplLeftThick = $ws.Range("ZZ100")   # will carry the "s=2" thick-left-border format
$tplLeftRed   = $ws.Range("ZZ101")   # will carry the "s=3" thin-red-left-border format

$ws.Range("E2").Copy()
$tplLeftThick.PasteSpecial(-4122)    # -4122 = xlPasteFormats

$ws.Range("C6").Copy()
$tplLeftRed.PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- 2. Clear the existing data block -------------------------------------
$ws.Range("A2:F17").Clear()

# --- 3. Rebuild the layout -------------------------------------------------
$cells = @(
    @{ Ref = 'B2'; Style = '2'; Text = '→' },
    @{ Ref = 'F2'; Style = ''; Text = 'CarPollutionPermit' },
    @{ Ref = 'B3'; Style = '2'; Text = '' },
    @{ Ref = 'F3'; Style = ''; Text = '__init__' },
    @{ Ref = 'B4'; Style = '2'; Text = '' },
    @{ Ref = 'F4'; Style = ''; Text = 'check_permit' },
    @{ Ref = 'B5'; Style = '2'; Text = '' },
    @{ Ref = 'C5'; Style = '3'; Text = '→' },
    @{ Ref = 'F5'; Style = ''; Text = 'BikePollutionPermit' },
    @{ Ref = 'B6'; Style = '2'; Text = '' },
    @{ Ref = 'C6'; Style = '3'; Text = '' },
    @{ Ref = 'F6'; Style = ''; Text = '__init__' },
    @{ Ref = 'B7'; Style = '2'; Text = '' },
    @{ Ref = 'C7'; Style = '3'; Text = '' },
    @{ Ref = 'F7'; Style = ''; Text = 'check_permit' },
    @{ Ref = 'B8'; Style = '2'; Text = '' },
    @{ Ref = 'C8'; Style = '3'; Text = '' },
    @{ Ref = 'E8'; Style = '3'; Text = '→' },
    @{ Ref = 'F8'; Style = ''; Text = 'TractorPollutionPermit' },
    @{ Ref = 'B9'; Style = '2'; Text = '' },
    @{ Ref = 'C9'; Style = '3'; Text = '' },
    @{ Ref = 'E9'; Style = '3'; Text = '' },
    @{ Ref = 'F9'; Style = ''; Text = 'fetch_tractor' },
    @{ Ref = 'B10'; Style = '2'; Text = '' },
    @{ Ref = 'C10'; Style = '3'; Text = '' },
    @{ Ref = 'E10'; Style = '3'; Text = '←' },
    @{ Ref = 'F10'; Style = ''; Text = 'TractorPesticides' },
    @{ Ref = 'B11'; Style = '2'; Text = '' },
    @{ Ref = 'C11'; Style = '3'; Text = '' },
    @{ Ref = 'F11'; Style = ''; Text = 'fetch_pesticides_permit' },
    @{ Ref = 'B12'; Style = '2'; Text = '' },
    @{ Ref = 'C12'; Style = '3'; Text = '→' },
    @{ Ref = 'D12'; Style = '3'; Text = '→' },
    @{ Ref = 'F12'; Style = ''; Text = 'Vehicle' },
    @{ Ref = 'B13'; Style = '2'; Text = '' },
    @{ Ref = 'C13'; Style = '3'; Text = '' },
    @{ Ref = 'D13'; Style = '3'; Text = '' },
    @{ Ref = 'F13'; Style = ''; Text = '__init__' },
    @{ Ref = 'B14'; Style = '2'; Text = '' },
    @{ Ref = 'C14'; Style = '3'; Text = '' },
    @{ Ref = 'D14'; Style = '3'; Text = '' },
    @{ Ref = 'F14'; Style = ''; Text = 'mileage_calculator' },
    @{ Ref = 'B15'; Style = '2'; Text = '←' },
    @{ Ref = 'C15'; Style = '3'; Text = '' },
    @{ Ref = 'D15'; Style = '3'; Text = '←' },
    @{ Ref = 'F15'; Style = ''; Text = 'Car' },
    @{ Ref = 'B16'; Style = '2'; Text = '' },
    @{ Ref = 'C16'; Style = '3'; Text = '' },
    @{ Ref = 'F16'; Style = ''; Text = '__init__' },
    @{ Ref = 'B17'; Style = '2'; Text = '' },
    @{ Ref = 'C17'; Style = '3'; Text = '' },
    @{ Ref = 'F17'; Style = ''; Text = 'pollution_permit' },
    @{ Ref = 'A18'; Style = '2'; Text = '→' },
    @{ Ref = 'B18'; Style = '2'; Text = '' },
    @{ Ref = 'C18'; Style = '3'; Text = '' },
    @{ Ref = 'F18'; Style = ''; Text = 'Farzi' },
    @{ Ref = 'A19'; Style = '2'; Text = '' },
    @{ Ref = 'B19'; Style = '2'; Text = '' },
    @{ Ref = 'C19'; Style = '3'; Text = '' },
    @{ Ref = 'F19'; Style = ''; Text = '__init__' },
    @{ Ref = 'A20'; Style = '2'; Text = '' },
    @{ Ref = 'B20'; Style = '2'; Text = '' },
    @{ Ref = 'C20'; Style = '3'; Text = '' },
    @{ Ref = 'F20'; Style = ''; Text = 'check_farzi' },
    @{ Ref = 'A21'; Style = '2'; Text = '←' },
    @{ Ref = 'B21'; Style = '2'; Text = '←' },
    @{ Ref = 'C21'; Style = '3'; Text = '←' },
    @{ Ref = 'F21'; Style = ''; Text = 'Bike' },
    @{ Ref = 'F22'; Style = ''; Text = '__init__' },
    @{ Ref = 'F23'; Style = ''; Text = 'pollution_permit' },
    @{ Ref = 'F24'; Style = ''; Text = 'check_farzi' }
)

foreach ($item in $cells) {
    $target = $ws.Range($item.Ref)
    if ($item.Style -eq '2') {
        $tplLeftThick.Copy()
        $target.PasteSpecial(-4122)
    } elseif ($item.Style -eq '3') {
        $tplLeftRed.Copy()
        $target.PasteSpecial(-4122)
    }
    if ($item.Text -ne '') {
        $target.Value = $item.Text
    }
}

$excel.CutCopyMode = 0

# --- 4. Remove scratch template cells --------------------------------------
$tplLeftThick.Clear()
$tplLeftRed.Clear()

Write-Host "Layout rebuilt."
